$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(3833, 3833, 3869, 3995, 3997, 4254, 4254, 4611, 4983, 4983, 4983, 4984, 4984, 5067)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
